# camparam.xlsx — add MOG/MOG2 background-subtractor option to the optical
# flow method parameter, plus the three new parameters it needs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Default ofmethod value changes from "farneback" to "mog2" (B2/C2 share
#    the same shared string, so both are updated together).
$ws.Range("B2:C2").Value = "mog2"

# 2) Extend the comment on A2 (the ofmethod help text) with the new options.
$cmt = $ws.Range("A2").Comment
$null = $cmt.Text($cmt.Text() + "`nmog: background subtractor (opencv2 only)`nmog2: background subtractor")

# 3) Add the three new parameter rows used by the MOG2 background subtractor.
$ws.Range("A26").Value = "nhistory"
$ws.Range("B26").Value = 100
$ws.Range("C26").Value = 100

$ws.Range("A27").Value = "nmixtures"
$ws.Range("B27").Value = 5
$ws.Range("C27").Value = 5

$ws.Range("A28").Value = "varThreshold"
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 1

# 4) Move the active selection to A30 (matches the saved cursor position).
$null = $ws.Range("A30").Select()
